$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the deleted "Suspended Sediment Concentration" row (old row 10);
# rows below it (old 11, 12) shift up to become new rows 10, 11.
$ws.Rows.Item(10).Delete()

# Row 2
$ws.Range("H2").Value = 0.962264150943396
$ws.Range("J2").Value = 8.41
$ws.Range("K2").Value = -0.149897400820793
$ws.Range("L2").Value = -0.27688710040959
$ws.Range("M2").Value = 0.107168689096305
$ws.Range("N2").Value = -1.78237099668006

# Row 3
$ws.Range("F3").Value = 0.0588096890395532
$ws.Range("H3").Value = 0.377358490566038
$ws.Range("J3").Value = 0.021
$ws.Range("K3").Value = 0.0008287013443288
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0.0020068681318681
$ws.Range("N3").Value = 3.94619687775657
$ws.Range("P3").Value = "Very unlikely improving"

# Row 4
$ws.Range("E4").Value = "ok"
$ws.Range("F4").Value = 0.0279650226665659
$ws.Range("G4").Value = 0.0377358490566038
$ws.Range("H4").Value = 0.773584905660377
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 120
$ws.Range("K4").Value = 20.7327758221434
$ws.Range("L4").Value = 3.30180824450497
$ws.Range("M4").Value = 46.7929133535552
$ws.Range("N4").Value = 17.2773131851195
$ws.Range("P4").Value = "Extremely unlikely improving"

# Row 5
$ws.Range("F5").Value = 0.0012794330537947
$ws.Range("G5").Value = 0.471698113207547
$ws.Range("H5").Value = 0.528301886792453
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 0.0113672845527755
$ws.Range("K5").Value = 0.0027162191484919
$ws.Range("M5").Value = 0.006274130593707
$ws.Range("N5").Value = 23.895057222164
$ws.Range("P5").Value = "Exceptionally unlikely improving"

# Row 6
$ws.Range("D6").Value = $true
$ws.Range("E6").Value = "ok"
$ws.Range("F6").Value = 0.0013607781251607
$ws.Range("G6").Value = 0.188679245283019
$ws.Range("H6").Value = 0.377358490566038
$ws.Range("K6").Value = 0.001003434065934
$ws.Range("L6").Value = 0.0005014329635701
$ws.Range("M6").Value = 0.0019730904361722
$ws.Range("N6").Value = 16.7239010989011
$ws.Range("P6").Value = "Exceptionally unlikely improving"

# Row 7
$ws.Range("E7").Value = "WARNING: Sen slope influenced by censored values"
$ws.Range("F7").Value = 0.0041536812099195
$ws.Range("J7").Value = 0.623
$ws.Range("K7").Value = 0.0696268198362147
$ws.Range("L7").Value = 0.0237949202063229
$ws.Range("M7").Value = 0.114032234404986
$ws.Range("N7").Value = 11.1760545483491
$ws.Range("P7").Value = "Exceptionally unlikely improving"

# Row 8
$ws.Range("D8").Value = $false
$ws.Range("F8").Value = 0.007800427951528
$ws.Range("H8").Value = 0.679245283018868
$ws.Range("J8").Value = 7.56
$ws.Range("K8").Value = -0.0441376299874621
$ws.Range("L8").Value = -0.0811782376535461
$ws.Range("M8").Value = -0.016297438812807
$ws.Range("N8").Value = -0.583831084490239

# Row 9
$ws.Range("F9").Value = 0.0010987966487159
$ws.Range("H9").Value = 0.867924528301887
$ws.Range("J9").Value = 0.73
$ws.Range("K9").Value = 0.0905578512396694
$ws.Range("L9").Value = 0.0443361749211696
$ws.Range("M9").Value = 0.144946551401756
$ws.Range("N9").Value = 12.4051851013246
$ws.Range("P9").Value = "Exceptionally unlikely improving"

# Row 10
$ws.Range("F10").Value = 0.0004471631267326
$ws.Range("H10").Value = 0.867924528301887
$ws.Range("J10").Value = 1.64
$ws.Range("K10").Value = 0.145698074277854
$ws.Range("L10").Value = 0.0465157781257331
$ws.Range("M10").Value = 0.190652472527473
$ws.Range("N10").Value = 8.88402891938135
$ws.Range("P10").Value = "Exceptionally unlikely improving"

# Row 11
$ws.Range("F11").Value = 0.0001597390115868
$ws.Range("H11").Value = 0.811320754716981
$ws.Range("J11").Value = 0.04
$ws.Range("K11").Value = 0.0131263282259011
$ws.Range("L11").Value = 0.0067948753429859
$ws.Range("M11").Value = 0.0201079830256992
$ws.Range("N11").Value = 32.8158205647528
$ws.Range("P11").Value = "Exceptionally unlikely improving"
